# Add a "success feedback" block of rows to the recommendation form sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(1582872, "Test One", "Yup"),
    @(1536237, "Test One", "Yup"),
    @(1519712, "Test One", "Oh yeah"),
    @(1565818, "Test One", "Oh yeah"),
    @(1506281, "Test One", "Oh yeah")
)

$startRow = 62
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
